$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Delivery_results")
$ws2 = $wb.Worksheets.Item("Total_emissions")

# --- Add two new header columns (F, G), copying the header format from an
# existing header cell (E1) so the new cells pick up the same bold/border
# style used by the rest of row 1 instead of creating a brand-new style. ---
$ws1.Range("E1").Copy() | Out-Null
$ws1.Range("F1:G1").PasteSpecial(-4122) | Out-Null

$ws1.Range("F1").Value = "Average_waiting_time_(minutes)"
$ws1.Range("G1").Value = "Average_queue_length"

# --- Update existing columns D (Number_of_parcels_delivered) and
# E (Percentage_of_parcels_delivered_(%)) with the new figures, and fill in
# the new F (Average_waiting_time_(minutes)) / G (Average_queue_length)
# columns for every data row. ---
$ws1.Range("D2").Value = 25
$ws1.Range("E2").Value = 40.32
$ws1.Range("F2").Value = 3.26
$ws1.Range("G2").Value = 8.6

$ws1.Range("D3").Value = 23
$ws1.Range("E3").Value = 38.33
$ws1.Range("F3").Value = 3.33
$ws1.Range("G3").Value = 7

$ws1.Range("D4").Value = 21
$ws1.Range("E4").Value = 48.84
$ws1.Range("F4").Value = 3.38
$ws1.Range("G4").Value = 9

$ws1.Range("D5").Value = 16
$ws1.Range("E5").Value = 88.89
$ws1.Range("F5").Value = 3.21
$ws1.Range("G5").Value = 5.5

$ws1.Range("D6").Value = 27
$ws1.Range("E6").Value = 72.97
$ws1.Range("F6").Value = 2.68
$ws1.Range("G6").Value = 7.37

# --- Total_emissions sheet: value stays materially the same (float
# precision only), re-assert it so the cell gets rewritten. ---
$ws2.Range("A2").Value = 988.19178327192742
